$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells in row 1 (columns D, E, F)
$ws.Range("D1").Value = "ORG_QUAR_IDENOLD"
$ws.Range("E1").Value = "ORG_QUAR_IDENNEW"
$ws.Range("F1").Value = "ORG_QUAR_STATUS"

# Match the left-aligned style used by the rest of the header row
$ws.Range("D1:F1").HorizontalAlignment = -4131

# Update the active selection to F2
$ws.Range("F2").Select()
